$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must be forced to Text
# format first so Excel stores them as the literal string (matching the
# source data, which keeps these as text cells) instead of silently
# converting to a numeric value.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.397.58"
$ws.Range("E2").Value = "  +3.25%  "
$ws.Range("D3").Value = "2.079.54"
$ws.Range("E3").Value = "  +3.42%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("D6").Value = "0.650"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").Value = "64.76"
$ws.Range("E7").Value = "  +2.60%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.395"
$ws.Range("D10").Value = "59.61"
$ws.Range("E10").Value = "  +0.75%  "
$ws.Range("E11").Value = "  +8.95%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "0.929"
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").Value = "23.43"
$ws.Range("E14").Value = "  +19.25%  "
$ws.Range("D15").Value = "14.84"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "2.384.32"
$ws.Range("E16").Value = "  +3.52%  "
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("D18").Value = "2.078.25"
$ws.Range("E18").Value = "  +3.42%  "
$ws.Range("D19").Value = "37.373.08"
$ws.Range("E19").Value = "  +3.39%  "
$ws.Range("D20").Value = "73.76"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  +5.27%  "
$ws.Range("E22").Value = "  +4.74%  "
$ws.Range("D23").Value = "240.07"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("D27").Value = "10.02"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("D28").Value = "20.81"
$ws.Range("E28").Value = "  +5.98%  "
$ws.Range("D29").Value = "161.87"
$ws.Range("E29").Value = "  -2.56%  "
$ws.Range("D30").Value = "0.128"
$ws.Range("E30").Value = "  +28.76%  "
$ws.Range("E31").Value = "  +2.21%  "
$ws.Range("D32").Value = "5.17"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("E33").Value = "  +3.61%  "
$ws.Range("D34").Value = "0.0628"
$ws.Range("E34").Value = "  +3.28%  "
$ws.Range("D35").Value = "4.66"
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("D36").Value = "2.58"
$ws.Range("E36").Value = "  +4.12%  "
$ws.Range("D37").Value = "6.46"
$ws.Range("E37").Value = "  +11.82%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "1.85"
$ws.Range("E38").Value = "  +2.84%  "
$ws.Range("B39").Value = "BinanceUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").Value = "3.04"
$ws.Range("E40").Value = "  +30.01%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.29"
$ws.Range("E41").Value = "  +4.02%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "0.102"
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("E44").Value = "  +5.19%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "17.35"
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0220"
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "8.01"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "96.25"
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("D49").Value = "1.404.82"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").Value = "2.95"
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "46.66"
$ws.Range("E51").Value = "  -1.05%  "
